# The "6c23b1b2-5b65-4d3e-847f-ca28e36d6cbc" image_tags entry (row 2) is
# being removed. The "index" column (A) keeps its original 1..10 sequence,
# while the uuid/name/tags columns (B:D) for every subsequent row move up
# by one row, and the now-duplicate trailing row disappears.
#
# Copy() (rather than Delete-with-shift) is used so each cell's value AND
# its formatting/style travel together, exactly like the row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the uuid / image_name / image_tags columns up by one row (B3:D12 -> B2:D11),
# leaving column A (the 1..10 index) untouched.
$ws.Range("B3:D12").Copy($ws.Range("B2:D11"))

# The old last row (12) is now a duplicate of row 11 - remove it so the
# table ends at row 11 again.
$ws.Rows(12).Delete()

# Match the cursor position left behind by the edit.
$ws.Range("B13").Select() | Out-Null
